# "Icon chart is working!" -- add 3 new rows (treatment_id 16) describing the
# new "calendarIcon" view type to the "Web Parameters" sheet, and move the
# selection the way Excel would leave it after typing them in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)   # "Web Parameters"
$ws.Activate()

# view_type / interaction / variable_amount labels used by the new rows
$viewType    = "calendarIcon"
$none        = "none"
$comment     = "Calendar year view with icon and no interaction."

# Row 29 -- treatment_id 16, position 1
$ws.Cells.Item(29, 1).Value  = 16
$ws.Cells.Item(29, 2).Value  = 1
$ws.Cells.Item(29, 3).Value  = $viewType
$ws.Cells.Item(29, 4).Value  = $none
$ws.Cells.Item(29, 5).Value  = $none
$ws.Cells.Item(29, 6).Value  = 300
$ws.Cells.Item(29, 8).Value  = 44593
$ws.Cells.Item(29, 9).Value  = 700
$ws.Cells.Item(29, 11).Value = 44703
$ws.Cells.Item(29, 12).Value = 1100
$ws.Cells.Item(29, 14).Value = 100
$ws.Cells.Item(29, 15).Value = 100
$ws.Cells.Item(29, 20).Value = 10
$ws.Cells.Item(29, 21).Value = 8
$ws.Cells.Item(29, 22).Value = $comment

# Row 30 -- treatment_id 16, position 2
$ws.Cells.Item(30, 1).Value  = 16
$ws.Cells.Item(30, 2).Value  = 2
$ws.Cells.Item(30, 3).Value  = $viewType
$ws.Cells.Item(30, 4).Value  = $none
$ws.Cells.Item(30, 5).Value  = $none
$ws.Cells.Item(30, 6).Value  = 500
$ws.Cells.Item(30, 8).Value  = 44621
$ws.Cells.Item(30, 9).Value  = 800
$ws.Cells.Item(30, 11).Value = 44724
$ws.Cells.Item(30, 12).Value = 1100
$ws.Cells.Item(30, 14).Value = 100
$ws.Cells.Item(30, 15).Value = 100
$ws.Cells.Item(30, 20).Value = 10
$ws.Cells.Item(30, 21).Value = 8
$ws.Cells.Item(30, 22).Value = $comment

# Row 31 -- treatment_id 16, position 3
$ws.Cells.Item(31, 1).Value  = 16
$ws.Cells.Item(31, 2).Value  = 3
$ws.Cells.Item(31, 3).Value  = $viewType
$ws.Cells.Item(31, 4).Value  = $none
$ws.Cells.Item(31, 5).Value  = $none
$ws.Cells.Item(31, 6).Value  = 300
$ws.Cells.Item(31, 8).Value  = 44652
$ws.Cells.Item(31, 9).Value  = 1000
$ws.Cells.Item(31, 11).Value = 44757
$ws.Cells.Item(31, 12).Value = 1100
$ws.Cells.Item(31, 14).Value = 100
$ws.Cells.Item(31, 15).Value = 100
$ws.Cells.Item(31, 20).Value = 10
$ws.Cells.Item(31, 21).Value = 8
$ws.Cells.Item(31, 22).Value = $comment

# Leave the view scrolled/selected the way it ended up after data entry:
# frozen header pane stays put, the bottom pane scrolls right so column V is
# visible, and the last-typed cell (V30) is the active selection.
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$ws.Range("V30").Select() | Out-Null
